$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "norm_outside_temperature"
$ws.Range("D2").Value = -12

[void]$ws.Range("D3").Select()
